# Update cryptocurrency price/volume figures (and fix the Aave/Algorand row
# ordering) per the scheduled "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the cells we touch so Excel does not reinterpret
# numeric-looking strings (e.g. "0.600", "46.179.83", "0.0847") as numbers,
# which would silently strip meaningful leading/trailing zeros or merge the
# thousand-separator dots used by this sheet's custom "Price" formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.179.83'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.606.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.61'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.56'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.600'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.19'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.11'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.18'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.997.28'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.595.28'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.923'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.90'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '46.301.56'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.80'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.90'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '289.73'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +14.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.19'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.48%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.54'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.93%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.06'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.89'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '39.38'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.22'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.30'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.62'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.20'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0843'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.20'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.80'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.123'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.57%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.04%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.99%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.62'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.05'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.20'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +11.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.110.71'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.54'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.49'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.11%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '109.34'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.203'
